$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format first so values such as
# "1.000" or "246.00" are stored verbatim instead of being parsed as numbers.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.540.84'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.919.57'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '246.00'
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.4795'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("D8").Value = '0.2905'
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("D9").Value = '0.06725'
$ws.Range("D10").Value = '110.75'
$ws.Range("E10").Value = '  +3.93%  '
$ws.Range("D11").Value = '19.06'
$ws.Range("E11").Value = '  +3.86%  '
$ws.Range("D12").Value = '1.915.13'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = '0.07571'
$ws.Range("E13").Value = '  -2.38%  '
$ws.Range("D14").Value = '5.283'
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").Value = '0.6680'
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").Value = '299.11'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("D17").Value = '30.509.88'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '5.648'
$ws.Range("E18").Value = '  +6.64%  '
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '0.000007583'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = '2.157.87'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '6.486'
$ws.Range("E24").Value = '  +4.29%  '
$ws.Range("D25").Value = '9.481'
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").Value = '164.82'
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("D27").Value = '20.32'
$ws.Range("E27").Value = '  -5.03%  '
$ws.Range("D28").Value = '2.114'
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").Value = '0.1078'
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = '1.402'
$ws.Range("E30").Value = '  +2.52%  '
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("D32").Value = '4.052'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").Value = '0.05004'
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("D34").Value = '0.7383'
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '2.732'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").Value = '0.02038'
$ws.Range("E38").Value = '  -2.78%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = '111.14'
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("D41").Value = '2.024'
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("D42").Value = '0.4440'
$ws.Range("E42").Value = '  +3.95%  '
$ws.Range("D43").Value = '72.53'
$ws.Range("E43").Value = '  +7.06%  '
$ws.Range("D44").Value = '0.8648'
$ws.Range("E44").Value = '  -1.06%  '
$ws.Range("D45").Value = '5.880'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").Value = '7.285'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").Value = '49.28'
$ws.Range("E48").Value = '  -1.76%  '
$ws.Range("D49").Value = '9.337'
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("D50").Value = '0.2569'
$ws.Range("E50").Value = '  +4.05%  '
$ws.Range("D51").Value = '0.1232'
$ws.Range("E51").Value = '  +0.84%  '

# Restore the cells to the workbook default (unstyled) look now that the
# text values are safely stored, matching the original formatting.
$fmtRange.Style = "Normal"

